$d = $word.ActiveDocument

# Remove the "Difficulty: Easy" paragraph entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Difficulty: Easy") {
        $p.Range.Delete()
        break
    }
}

# Add _GoBack bookmark at start of "Preview: True" paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Preview: True") {
        $r = $p.Range.Duplicate
        $r.Collapse(1)
        $d.Bookmarks.Add("_GoBack", $r)
        break
    }
}

Write-Host "Paragraphs count: " $d.Paragraphs.Count
foreach ($p in $d.Paragraphs) {
    Write-Host "----"
    Write-Host $p.Range.Text
}
